$wb = $excel.ActiveWorkbook

$wsFileLocations = $wb.Worksheets.Item("FileLocations")
$wsSite           = $wb.Worksheets.Item("Site")
$wsVendorGeneral  = $wb.Worksheets.Item("Vendor General")
$wsValues         = $wb.Worksheets.Item("Values")
$wsGamesa         = $wb.Worksheets.Item("Gamesa")

# --- Content changes -------------------------------------------------

# FileLocations!B4 - directory path now points at the Gamesa subfolder
$wsFileLocations.Range("B4").Value = "C:\\Users\\vsivakumaran\\Documents\\Test c#\\Gamesa"

# FileLocations!B column got wider to fit the longer path
$wsFileLocations.Columns.Item(2).ColumnWidth = 51.8

# Values!A6 - "Regulator" renamed to "Regulatory"
$wsValues.Range("A6").Value = "Regulatory"

# --- View / selection changes -----------------------------------------

# FileLocations is no longer the active tab; last selection left on B26
$null = $wsFileLocations.Range("B26").Select()

# Vendor General's last selection left on E26
$null = $wsVendorGeneral.Range("E26").Select()

# Values sheet now has a remembered selection at C13
$null = $wsValues.Range("C13").Select()

# Gamesa's last selection moved from F4 to F5
$null = $wsGamesa.Range("F5").Select()

# Site becomes the active tab, with selection left on A3 (select last so
# it "wins" the workbook-level active tab / tabSelected flag)
$null = $wsSite.Range("A3").Select()
